# CH-140 Golden Period.xlsx — "Understanding" sheet update
#
# Adds two more experiment blocks to the bottom of the sheet showing that
# INDEX(x,1,1) and SINGLE(+x) can both stand in for the old (deprecated)
# bare SINGLE usage inside GROUPBY's aggregation LAMBDA, plus a short note
# ("So INDEX can replace single") above the first of the two blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Understanding")

# ---------------------------------------------------------------------
# Row 68: short note introducing the next experiment
# ---------------------------------------------------------------------
$ws.Range("L68").Value = "So INDEX can replace single"

# ---------------------------------------------------------------------
# Block 1 (rows 69-80): GROUPBY(...,LAMBDA(x,TEXTJOIN(",",,INDEX(x,1,1))))
# ---------------------------------------------------------------------
$ws.Range("F69").Value = "A"
$ws.Range("G69").Value = "z"
$ws.Range("H69").Value = 11

$ws.Range("F70").Value = "A"
$ws.Range("G70").Formula = '=G69&"Z"'
$ws.Range("H70").Value = 10

$ws.Range("F71").Value = "A"
$ws.Range("G71").Formula = '=G70&"Z"'
$ws.Range("H71").Value = 9

$ws.Range("F72").Value = "B"
$ws.Range("G72").Formula = '=G71&"Z"'
$ws.Range("H72").Value = 8

$ws.Range("F73").Value = "B"
$ws.Range("G73").Formula = '=G72&"Z"'
$ws.Range("H73").Value = 7

$ws.Range("F74").Value = "B"
$ws.Range("G74").Formula = '=G73&"Z"'
$ws.Range("H74").Value = 6

$ws.Range("F75").Value = "C"
$ws.Range("G75").Formula = '=G74&"Z"'
$ws.Range("H75").Value = 5

$ws.Range("F76").Value = "C"
$ws.Range("G76").Formula = '=G75&"Z"'
$ws.Range("H76").Value = 4

$ws.Range("F77").Value = "D"
$ws.Range("G77").Formula = '=G76&"Z"'
$ws.Range("H77").Value = 3

$ws.Range("F78").Value = "D"
$ws.Range("G78").Formula = '=G77&"Z"'
$ws.Range("H78").Value = 2

$ws.Range("F79").Value = "D"
$ws.Range("G79").Formula = '=G78&"Z"'
$ws.Range("H79").Value = 1

$ws.Range("L69:N73").FormulaArray = '=GROUPBY(F69:F79,G69:H79,LAMBDA(x,TEXTJOIN(",",,INDEX(x,1,1))))'

# Legend-style notes to the right of the spilled GROUPBY result
$ws.Range("M76").Value = "z,zZ,zZZ"
$ws.Range("M77").Value = "zZZZ,zZZZZ,zZZZZZ"
$ws.Range("M78").Value = "zZZZZZZ,zZZZZZZZ"
$ws.Range("M79").Value = "zZZZZZZZZ,zZZZZZZZZZ,zZZZZZZZZZZ"
$ws.Range("M80").Value = "z,zZ,zZZ,zZZZ,zZZZZ,zZZZZZ,zZZZZZZ,zZZZZZZZ,zZZZZZZZZ,zZZZZZZZZZ,zZZZZZZZZZZ"

# ---------------------------------------------------------------------
# Row 82: shows the actual formula text used in the block below
# (mirrors how Excel renders SINGLE(+x) back as the @ operator)
# ---------------------------------------------------------------------
$ws.Range("L82").Formula = '=FORMULATEXT(L83)'

# ---------------------------------------------------------------------
# Block 2 (rows 83-94): GROUPBY(...,LAMBDA(x,TEXTJOIN(",",,SINGLE(+x))))
# ---------------------------------------------------------------------
$ws.Range("F83").Value = "A"
$ws.Range("G83").Value = "z"
$ws.Range("H83").Value = 11

$ws.Range("F84").Value = "A"
$ws.Range("G84").Formula = '=G83&"Z"'
$ws.Range("H84").Value = 10

$ws.Range("F85").Value = "A"
$ws.Range("G85").Formula = '=G84&"Z"'
$ws.Range("H85").Value = 9

$ws.Range("F86").Value = "B"
$ws.Range("G86").Formula = '=G85&"Z"'
$ws.Range("H86").Value = 8

$ws.Range("F87").Value = "B"
$ws.Range("G87").Formula = '=G86&"Z"'
$ws.Range("H87").Value = 7

$ws.Range("F88").Value = "B"
$ws.Range("G88").Formula = '=G87&"Z"'
$ws.Range("H88").Value = 6

$ws.Range("F89").Value = "C"
$ws.Range("G89").Formula = '=G88&"Z"'
$ws.Range("H89").Value = 5

$ws.Range("F90").Value = "C"
$ws.Range("G90").Formula = '=G89&"Z"'
$ws.Range("H90").Value = 4

$ws.Range("F91").Value = "D"
$ws.Range("G91").Formula = '=G90&"Z"'
$ws.Range("H91").Value = 3

$ws.Range("F92").Value = "D"
$ws.Range("G92").Formula = '=G91&"Z"'
$ws.Range("H92").Value = 2

$ws.Range("F93").Value = "D"
$ws.Range("G93").Formula = '=G92&"Z"'
$ws.Range("H93").Value = 1

$ws.Range("L83:N87").FormulaArray = '=GROUPBY(F83:F93,G83:H93,LAMBDA(x,TEXTJOIN(",",,SINGLE(+x))))'

# Legend-style notes to the right of the spilled GROUPBY result
$ws.Range("M90").Value = "z,zZ,zZZ"
$ws.Range("M91").Value = "zZZZ,zZZZZ,zZZZZZ"
$ws.Range("M92").Value = "zZZZZZZ,zZZZZZZZ"
$ws.Range("M93").Value = "zZZZZZZZZ,zZZZZZZZZZ,zZZZZZZZZZZ"
$ws.Range("M94").Value = "z,zZ,zZZ,zZZZ,zZZZZ,zZZZZZ,zZZZZZZ,zZZZZZZZ,zZZZZZZZZ,zZZZZZZZZZ,zZZZZZZZZZZ"

# ---------------------------------------------------------------------
# View state: scroll down to the newly added content and select R76
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("R76").Select()
